$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 1484
$wsExhibit.Range("F4").Value = 2113
$wsExhibit.Range("F5").Value = 7409
$wsExhibit.Range("F7").Value = 4768
$wsExhibit.Range("F8").Value = 7015
$wsExhibit.Range("F10").Value = 273
$wsExhibit.Range("F11").Value = 1487
$wsExhibit.Range("F12").Value = 857
$wsExhibit.Range("F13").Value = 169
$wsExhibit.Range("F14").Value = 52
$wsExhibit.Range("F15").Value = 1162
$wsExhibit.Range("F17").Value = 159
$wsExhibit.Range("F18").Value = 6
$wsExhibit.Range("F21").Value = 1156
$wsExhibit.Range("F22").Value = 753
$wsExhibit.Range("F24").Value = 46
$wsExhibit.Range("F25").Value = 1225
$wsExhibit.Range("F30").Value = 172
$wsExhibit.Range("F32").Value = 35
$wsExhibit.Range("F33").Value = 91
$wsExhibit.Range("F34").Value = 32
$wsExhibit.Range("F35").Value = 546
$wsExhibit.Range("F39").Value = 370
$wsExhibit.Range("F40").Value = 1198
$wsExhibit.Range("F41").Value = 575
$wsExhibit.Range("F42").Value = 138
$wsExhibit.Range("F44").Value = 19

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F3").Value = 33
$wsShow.Range("F33").Value = 608
$wsShow.Range("G41").Value = 171

$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F8").Value = 1575
$wsLocal.Range("F9").Value = 2472

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 1484
$wsAll.Range("F9").Value = 7409
$wsAll.Range("F11").Value = 4768
$wsAll.Range("F13").Value = 7015
$wsAll.Range("F14").Value = 273
$wsAll.Range("F15").Value = 1487
$wsAll.Range("F16").Value = 857
$wsAll.Range("F17").Value = 169
$wsAll.Range("F18").Value = 1575
$wsAll.Range("F19").Value = 2472
$wsAll.Range("F21").Value = 52
$wsAll.Range("F22").Value = 1162
$wsAll.Range("F23").Value = 159
$wsAll.Range("F24").Value = 6
$wsAll.Range("F25").Value = 224
$wsAll.Range("F26").Value = 1156
$wsAll.Range("F28").Value = 753
$wsAll.Range("F30").Value = 1225
$wsAll.Range("F32").Value = 172
$wsAll.Range("F35").Value = 35
$wsAll.Range("F36").Value = 91
$wsAll.Range("F38").Value = 546
$wsAll.Range("F39").Value = 608
$wsAll.Range("F43").Value = 370
$wsAll.Range("F44").Value = 575
$wsAll.Range("G47").Value = 171
$wsAll.Range("F48").Value = 138
